# Updates cryptos.xlsx "Price" (D) and "Volume(1h)" (E) columns with refreshed
# coinranking.com figures, and re-ranks rows 43-45 (FTXToken / VeChain / Maker)
# to reflect the new coin order. Matches the "Updated cryptos list" GH Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many "Price" values look like plain numbers/dates to Excel (e.g. "1.00", "0.625").
# The sheet stores them as text, so force text storage via NumberFormat "@" before
# assigning any Price cell whose new value Excel would otherwise auto-convert.

# Row 2
$ws.Range("D2").Value = "37.815.22"
$ws.Range("E2").Value = "  +1.14%  "
# Row 3
$ws.Range("D3").Value = "2.084.82"
# Row 4
$ws.Range("E4").Value = "  -0.04%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.25"
$ws.Range("E5").Value = "  -0.42%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  -0.28%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.10"
$ws.Range("E7").Value = "  +2.83%  "
# Row 8
$ws.Range("E8").Value = "  -0.06%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  -0.40%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0790"
$ws.Range("E10").Value = "  +2.17%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +2.77%  "
# Row 12
$ws.Range("D12").Value = "2.390.89"
$ws.Range("E12").Value = "  +0.75%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.70"
$ws.Range("E13").Value = "  +1.73%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.24"
$ws.Range("E14").Value = "  +1.99%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.772"
$ws.Range("E15").Value = "  -0.67%  "
# Row 16
$ws.Range("E16").Value = "  +2.29%  "
# Row 17
$ws.Range("D17").Value = "2.075.54"
$ws.Range("E17").Value = "  +0.39%  "
# Row 18
$ws.Range("D18").Value = "37.726.81"
$ws.Range("E18").Value = "  +1.04%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("E19").Value = "  -0.03%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.40"
$ws.Range("E20").Value = "  +2.62%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("E21").Value = "  +1.50%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.98"
$ws.Range("E22").Value = "  +0.87%  "
# Row 23
$ws.Range("E23").Value = "  -0.03%  "
# Row 24
$ws.Range("E24").Value = "  -0.76%  "
# Row 25
$ws.Range("E25").Value = "  -1.15%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.44"
$ws.Range("E26").Value = "  +1.91%  "
# Row 27
$ws.Range("E27").Value = "  +8.41%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.03"
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.41"
$ws.Range("E29").Value = "  +0.31%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.55"
$ws.Range("E30").Value = "  +2.29%  "
# Row 31
$ws.Range("E31").Value = "  +2.41%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.71"
$ws.Range("E32").Value = "  +3.61%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.74"
$ws.Range("E33").Value = "  +4.31%  "
# Row 34
$ws.Range("E34").Value = "  +2.03%  "
# Row 35
$ws.Range("E35").Value = "  +1.94%  "
# Row 36
$ws.Range("E36").Value = "  +3.30%  "
# Row 37
$ws.Range("E37").Value = "  +2.94%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.06%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  -3.43%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0992"
$ws.Range("E40").Value = "  +3.58%  "
# Row 41
$ws.Range("E41").Value = "  -0.11%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.72"
$ws.Range("E42").Value = "  +1.03%  "
# Row 43
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.36"
$ws.Range("E43").Value = "  +4.14%  "
# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0215"
$ws.Range("E44").Value = "  +1.02%  "
# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.459.83"
$ws.Range("E45").Value = "  -1.78%  "
# Row 46
$ws.Range("E46").Value = "  +1.15%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.08"
$ws.Range("E47").Value = "  +6.35%  "
# Row 48
$ws.Range("E48").Value = "  +4.21%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.41"
$ws.Range("E49").Value = "  +2.90%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.04"
$ws.Range("E50").Value = "  +2.66%  "
# Row 51
$ws.Range("D51").Value = "2.275.48"
$ws.Range("E51").Value = "  +0.64%  "
